# Refresh the cryptocurrency price/volume snapshot in the "cryptos" worksheet.
# For Price-column (D) values that look like plain numbers (e.g. "218.68"), a
# leading apostrophe is used so Excel stores them as literal text (matching the
# source data, which mixes thousands-separated strings like "26.161.12" with
# plain decimals) instead of silently re-typing them as numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.161.12"
$ws.Range("E2").Value = "  -0.01%  "

# Row 3
$ws.Range("D3").Value = "1.655.99"
$ws.Range("E3").Value = "  -0.17%  "

# Row 4
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").Value = "'218.68"
$ws.Range("E5").Value = "  -0.30%  "

# Row 6
$ws.Range("D6").Value = "'0.5235"
$ws.Range("E6").Value = "  +0.18%  "

# Row 7
$ws.Range("E7").Value = "  -0.22%  "

# Row 8
$ws.Range("D8").Value = "'0.2665"
$ws.Range("E8").Value = "  +1.29%  "

# Row 9
$ws.Range("D9").Value = "'0.06343"
$ws.Range("E9").Value = "  +0.73%  "

# Row 10
$ws.Range("E10").Value = "  -0.45%  "

# Row 11
$ws.Range("D11").Value = "'0.07678"
$ws.Range("E11").Value = "  -1.77%  "

# Row 12
$ws.Range("D12").Value = "'4.618"
$ws.Range("E12").Value = "  +2.73%  "

# Row 13
$ws.Range("D13").Value = "1.710.55"
$ws.Range("E13").Value = "  +3.27%  "

# Row 14
$ws.Range("D14").Value = "1.884.37"
$ws.Range("E14").Value = "  -0.13%  "

# Row 15
$ws.Range("D15").Value = "'0.5610"
$ws.Range("E15").Value = "  +1.09%  "

# Row 16
$ws.Range("D16").Value = "0.0₅8184"
$ws.Range("E16").Value = "  +2.13%  "

# Row 17
$ws.Range("D17").Value = "'65.52"
$ws.Range("E17").Value = "  +0.62%  "

# Row 18
$ws.Range("D18").Value = "26.151.58"
$ws.Range("E18").Value = "  -0.14%  "

# Row 19
$ws.Range("E19").Value = "  -0.25%  "

# Row 20
$ws.Range("D20").Value = "'4.651"
$ws.Range("E20").Value = "  +0.25%  "

# Row 21
$ws.Range("D21").Value = "'10.47"
$ws.Range("E21").Value = "  +3.55%  "

# Row 22
$ws.Range("D22").Value = "'192.96"
$ws.Range("E22").Value = "  -1.55%  "

# Row 23
$ws.Range("D23").Value = "'5.955"
$ws.Range("E23").Value = "  -0.10%  "

# Row 24
$ws.Range("E24").Value = "  -0.23%  "

# Row 25
$ws.Range("D25").Value = "'145.31"
$ws.Range("E25").Value = "  -0.71%  "

# Row 26
$ws.Range("E26").Value = "  -0.59%  "

# Row 27
$ws.Range("D27").Value = "'7.261"
$ws.Range("E27").Value = "  +1.66%  "

# Row 28
$ws.Range("D28").Value = "'15.95"
$ws.Range("E28").Value = "  -0.27%  "

# Row 29
$ws.Range("D29").Value = "'1.518"
$ws.Range("E29").Value = "  +1.69%  "

# Row 30
$ws.Range("D30").Value = "'0.05487"
$ws.Range("E30").Value = "  -4.39%  "

# Row 31
$ws.Range("D31").Value = "'1.271"
$ws.Range("E31").Value = "  -0.34%  "

# Row 32
$ws.Range("D32").Value = "'3.470"
$ws.Range("E32").Value = "  -0.72%  "

# Row 33
$ws.Range("E33").Value = "  -0.09%  "

# Row 34
$ws.Range("D34").Value = "'1.564"

# Row 35
$ws.Range("D35").Value = "'0.9497"
$ws.Range("E35").Value = "  -0.62%  "

# Row 36
$ws.Range("D36").Value = "'2.780"
$ws.Range("E36").Value = "  -0.92%  "

# Row 37
$ws.Range("E37").Value = "  -0.71%  "

# Row 38
$ws.Range("D38").Value = "'0.5683"
$ws.Range("E38").Value = "  -0.54%  "

# Row 39
$ws.Range("D39").Value = "'0.01589"
$ws.Range("E39").Value = "  -0.43%  "

# Row 40
$ws.Range("E40").Value = "  -1.51%  "

# Row 41
$ws.Range("E41").Value = "  -0.22%  "

# Row 42
$ws.Range("D42").Value = "1.028.71"
$ws.Range("E42").Value = "  -3.40%  "

# Row 43
$ws.Range("D43").Value = "'0.8306"
$ws.Range("E43").Value = "  -2.28%  "

# Row 44
$ws.Range("D44").Value = "'100.96"

# Row 45
$ws.Range("D45").Value = "1.794.79"
$ws.Range("E45").Value = "  -0.14%  "

# Row 46
$ws.Range("D46").Value = "'57.87"
$ws.Range("E46").Value = "  -0.30%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈104"
$ws.Range("E47").Value = "  +1.28%  "

# Row 48
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'0.9972"
$ws.Range("E48").Value = "  -0.83%  "

# Row 49
$ws.Range("D49").Value = "'0.4347"
$ws.Range("E49").Value = "  -1.32%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.012"
$ws.Range("E50").Value = "  -0.35%  "

# Row 51
$ws.Range("D51").Value = "'0.05208"
$ws.Range("E51").Value = "  +0.14%  "
